# ============================================================================
# feat: add 2022-Q4 data
#
# Inserts a new worksheet "2022-Q4" right after "总计" (pushing the existing
# "2022-Q2" / "2022-Q1" sheets one slot to the right), fills it with the
# Q4 fund-holding data, and updates the "总计" (totals) sheet with a new
# row for 2022-Q4 (inserted above the existing 2022-Q2 / 2022-Q1 rows).
#
# NOTE: worksheet object references returned by Worksheets.Item(...) are
# resolved live against the *current* sheet order, so any reference fetched
# before a sheet-insertion can silently start pointing at a different sheet
# once the collection shifts. To stay safe we do the insert/rename first,
# with nothing else in between, and only *afterwards* fetch the stable
# by-name references used for the rest of the script.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
#    current "2022-Q2" sheet), matching the target tab order:
#      总计, 2022-Q4, 2022-Q2, 2022-Q1
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Re-fetch stable references now that the sheet order/names are final.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q4    = $wb.Worksheets.Item("2022-Q4")
$q2    = $wb.Worksheets.Item("2022-Q2")

# Clone number/border/font formatting from the "2022-Q2" sheet (style used
# for header row B1:H1 and the leading-index column A2:A4) so the new sheet
# matches the look of its siblings.
$q2.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q2.Range("A2:A3").Copy()
$q4.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Populate "2022-Q4" sheet contents.
# ---------------------------------------------------------------------------
# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2 - 006693 金信消费升级股票C
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "006693"
$q4.Range("C2").Value = "金信消费升级股票C"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.97"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "87.56"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "5.44"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0528"
$q4.Range("H2").Value = 9

# Row 3 - 006692 金信消费升级股票A
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "006692"
$q4.Range("C3").Value = "金信消费升级股票A"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.75"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "87.56"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "5.44"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0408"
$q4.Range("H3").Value = 9

# Row 4 - 002862 金信量化精选灵活配置混合
$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "002862"
$q4.Range("C4").Value = "金信量化精选灵活配置混合"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.46"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "80.41"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "3.38"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.0155"
$q4.Range("H4").Value = 8

# ---------------------------------------------------------------------------
# 4) Update "总计" sheet: shift the existing 2022-Q2 / 2022-Q1 rows down by
#    one and insert the new 2022-Q4 totals row above them.
#    Before:                       After:
#      row2 = 2022-Q2, 2, 0.04       row2 = 2022-Q4, 3, 0.11
#      row3 = 2022-Q1, 2, 0.04       row3 = 2022-Q2, 2, 0.04
#                                     row4 = 2022-Q1, 2, 0.04
# ---------------------------------------------------------------------------
# Copy the index-column style (A2, already s="2") down onto the newly used A4
# before overwriting values, so every row in A keeps the same look.
$total.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.04

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.04

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.11
